$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "ww"
$ws.Range("C8").Value = "ww"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
